# Rows 2-9 (the data rows; row 1 is the header) have had their entire
# record contents reassigned to different row positions -- i.e. the 8
# observation records got permuted across rows 2..9. The mapping below is
# target-row <- source-row (source-row = which row currently, before this
# edit, holds the data that must end up at target-row), derived by matching
# each record's unique Id (column A):
#
#   2 <- 4      3 <- 7      4 <- 3      5 <- 9
#   6 <- 5      7 <- 8      8 <- 2      9 <- 6
#
# Only these columns ever contain data in rows 2-9 (verified against the
# sheet): A B C D E F G H I M P Q R S T U V W Y Z AA AB AD AE AG AT AW AX AY
# (every other column between A and AY, e.g. J/K/L/N/O/AC/AF/..., is never
# populated in any data row and must stay that way.)
#
# Because this permutation has two cycles ((2 4 3 7 8) and (5 9 6)), an
# in-place cell-by-cell copy would clobber source rows before they get used.
# So: stage every source row's relevant cells into scratch rows first
# (single-cell Copy preserves the exact stored type/format, unlike reading
# .Value() into an array and writing it back, which lets Excel "smart"
# re-interpret text dates as date serials and drops empty cells), then copy
# from the scratch rows into the final destinations, then clear the scratch
# rows again.
#
# Column M ("Aktivitet") is sparse -- only rows with an activity note carry
# a cell there at all; the rest have no <c> element for M whatsoever (not
# even an empty one). Cells.Item(r,c).Value() reads back $null for a truly
# absent cell vs. "" for a present-but-empty one, so that distinguishes the
# two cases. A .Copy() of an absent cell onto another cell is a no-op (it
# does NOT clear the destination) -- but copying it onto a previously-absent
# scratch/destination cell would still materialize an empty stub there, so
# absent source cells are simply skipped everywhere below instead of copied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataCols = @(1,2,3,4,5,6,7,8,9,13,16,17,18,19,20,21,22,23,25,26,27,28,30,31,33,46,49,50,51)
# A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 I=9 M=13 P=16 Q=17 R=18 S=19 T=20 U=21 V=22
# W=23 Y=25 Z=26 AA=27 AB=28 AD=30 AE=31 AG=33 AT=46 AW=49 AX=50 AY=51

$mapping = @{
    2 = 4
    3 = 7
    4 = 3
    5 = 9
    6 = 5
    7 = 8
    8 = 2
    9 = 6
}

$scratchOffset = 100   # source row R stages at row R + 100 (well clear of used data)

# 1) Stage each distinct source row's data cells into its scratch row.
$sourceRows = $mapping.Values | Sort-Object -Unique
foreach ($srcRow in $sourceRows) {
    $scratchRow = $srcRow + $scratchOffset
    foreach ($col in $dataCols) {
        $srcCell = $ws.Cells.Item($srcRow, $col)
        if ($srcCell.Value() -ne $null) {
            $srcCell.Copy($ws.Cells.Item($scratchRow, $col))
        }
    }
}

# 2) A handful of cells are present (non-empty) today but must end up
#    completely absent, because the record now landing on that row never
#    had that column populated. Copying a genuinely-absent source cell is a
#    no-op on the destination (it does NOT clear existing content), so those
#    specific cells need an explicit clear first. Every other cell either
#    keeps the same presence/absence it already has, or gets overwritten by
#    the copy step below -- nothing else needs this.
$ws.Cells.Item(4, 13).ClearContents()   # M4
$ws.Cells.Item(5, 13).ClearContents()   # M5

# 3) Copy staged data into the final destination rows (skipping cells that
#    were never staged because the source was genuinely absent).
foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $scratchRow = $srcRow + $scratchOffset
    foreach ($col in $dataCols) {
        $scratchCell = $ws.Cells.Item($scratchRow, $col)
        if ($scratchCell.Value() -ne $null) {
            $scratchCell.Copy($ws.Cells.Item($targetRow, $col))
        }
    }
}

# 4) Remove the scratch rows' contents so they don't linger in the sheet.
#    (Deleting whole rows would shift everything below up on each call --
#    order-dependent and easy to get wrong -- so just clear the cells we
#    actually wrote instead; an empty row is indistinguishable from an
#    absent one for xlsx's sparse row storage.)
foreach ($srcRow in $sourceRows) {
    $scratchRow = $srcRow + $scratchOffset
    foreach ($col in $dataCols) {
        $ws.Cells.Item($scratchRow, $col).ClearContents()
    }
}
